$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value2 = "63.114.95"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value2 = "  -1.49%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value2 = "3.253.70"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value2 = "  +3.29%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value2 = "0.999"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value2 = "  -0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value2 = "593.72"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value2 = "  -1.85%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value2 = "140.72"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value2 = "  -2.05%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value2 = "  +0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value2 = "3.247.78"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value2 = "  +3.26%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value2 = "0.520"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value2 = "  -0.93%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value2 = "0.148"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value2 = "  -1.15%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value2 = "5.37"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value2 = "  -0.47%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value2 = "0.464"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value2 = "  -0.52%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value2 = "0.0000249"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value2 = "  -2.21%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value2 = "34.50"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value2 = "  -1.50%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value2 = "3.778.78"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value2 = "  +2.92%  "
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value2 = "  +0.10%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value2 = "3.244.28"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value2 = "  +3.28%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value2 = "63.137.28"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value2 = "  -1.52%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value2 = "6.76"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value2 = "  -1.10%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value2 = "476.08"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value2 = "  -2.65%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value2 = "14.15"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value2 = "  -3.70%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value2 = "0.722"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value2 = "  +1.44%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value2 = "7.92"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value2 = "  +3.63%  "
$ws.Range("B24").Value2 = "Litecoin"
$ws.Range("C24").Value2 = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value2 = "83.72"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value2 = "  -4.43%  "
$ws.Range("B25").Value2 = "InternetComputer(DFINITY)"
$ws.Range("C25").Value2 = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value2 = "13.30"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value2 = "  -0.16%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value2 = "  -0.02%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value2 = "2.73"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value2 = "  -1.25%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value2 = "7.30"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value2 = "  +4.32%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value2 = "8.09"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value2 = "  -1.46%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value2 = "2.13"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value2 = "  +3.05%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value2 = "27.59"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value2 = "  +0.73%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value2 = "1.00"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value2 = "  -0.08%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value2 = "0.108"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value2 = "  -2.55%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value2 = "2.54"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value2 = "  -3.68%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value2 = "1.09"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value2 = "  -1.15%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value2 = "5.88"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value2 = "  -2.40%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value2 = "52.79"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value2 = "  +0.16%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value2 = "0.0₃0713"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value2 = "  -4.52%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value2 = "0.0392"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value2 = "  -1.39%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value2 = "420.16"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value2 = "  -3.89%  "
$ws.Range("B41").Value2 = "Maker"
$ws.Range("C41").Value2 = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value2 = "2.991.35"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value2 = "  +2.28%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value2 = "8.39"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value2 = "  +1.00%  "
$ws.Range("B43").Value2 = "dogwifhat"
$ws.Range("C43").Value2 = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value2 = "2.74"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value2 = "  -7.73%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value2 = "0.112"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value2 = "  -6.89%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value2 = "0.268"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value2 = "  +3.32%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value2 = "2.16"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value2 = "  -1.79%  "
$ws.Range("B47").Value2 = "USDe"
$ws.Range("C47").Value2 = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value2 = "0.999"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value2 = "  +0.04%  "
$ws.Range("B48").Value2 = "InjectiveProtocol"
$ws.Range("C48").Value2 = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value2 = "25.91"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value2 = "  -0.25%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value2 = "2.31"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value2 = "  -4.36%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value2 = "0.114"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value2 = "  +0.05%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value2 = "119.10"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value2 = "  -0.92%  "
